$wb = $excel.ActiveWorkbook

# --- Commodity sheet: move selection from D4 to E4 ---
$wsCommodity = $wb.Worksheets.Item("Commodity")
$null = $wsCommodity.Range("E4").Select()

# --- Storage sheet: set selection to B37 ---
$wsStorage = $wb.Worksheets.Item("Storage")
$null = $wsStorage.Range("B37").Select()

# --- Process sheet: fix/clarify a few comments and correct bug in row 5 ---
$wsProcess = $wb.Worksheets.Item("Process")

# Bug fix: "Elec heating plant" (row 5) had zero fixed/specific investment cost values
$wsProcess.Range("B5").Value = 5000
$wsProcess.Range("C5").Value = 0.1

# Clarify which comments refer to the Richter source
$wsProcess.Range("H5").Value = "not mentioned in Richter"
$wsProcess.Range("H6").Value = "not mentioned in Richter (interesting: why?)"
$wsProcess.Range("H7").Value = "not mentioned in Richter"
$wsProcess.Range("H8").Value = "in Richter: BK"
$wsProcess.Range("H9").Value = "in Richter: MV"

# Process becomes the active sheet/tab, with selection on B6
$null = $wsProcess.Activate()
$null = $wsProcess.Range("B6").Select()
